# Updated symbol list on Sat Jan  7 10:59:46 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) hold numeric-looking text ("261.01", "1.82%", ...).
# The source workbook stores these as plain text (inlineStr), not numbers, so we prefix
# with a leading apostrophe to force Excel to keep them as text, then reset the cell
# style back to "Normal" to drop the quotePrefix flag that the apostrophe entry adds -
# this keeps the cell text-typed while matching the original (unstyled) appearance.

$ws.Range("D2").Value = "'261.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.82%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'2.04%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.707"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.79%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06093"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.28%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.671"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.47%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'3.149"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.39%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.8454"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.67%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").Value = "'0.9279"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.39%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1406"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.85%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.04851"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'15.83%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07105"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.42%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03075"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.72%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09066"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.48%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001543"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.36%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0006074"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.03%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006022"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.38%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.451"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.61%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.182"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.75%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'2.24%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D22").Value = "'4.099"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.32%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04254"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.13%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001223"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.34%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.003798"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'5.00%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'0.00%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'3.36%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03872"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.39%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1112"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.27%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.004083"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-35.12%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.01626"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'15.14%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'0.73%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005150"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-3.49%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.1355"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-43.90%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'23.71%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("E50").Style = "Normal"
